# extracting data from excel for a specific test case
#
# Header row gains "password"/"gender"/"DOB" columns (replacing the old
# "lastName" column, which is folded away), the old "Test3" sample row
# (g/h/i) is replaced by a real filled-in form submission row, and the
# selection / used-range grow from D:7 to F:7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header -----------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "firstName"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "gender"
$ws.Range("F1").Value = "DOB"

# --- Row 2: unchanged ----------------------------------------------------
$ws.Range("A2").Value = "Test1"
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "b"
$ws.Range("D2").Value = "c"

# --- Row 3: unchanged ----------------------------------------------------
$ws.Range("A3").Value = "Test2"
$ws.Range("B3").Value = "d"
$ws.Range("C3").Value = "e"
$ws.Range("D3").Value = "f"

# --- Row 4: new filled-in test data row (replaces old "Test3" g/h/i row) -
$ws.Range("A4").Value = "test_HomePage_FillingForm"
$ws.Range("B4").Value = "Abdur Rehman"
$ws.Range("C4").Value = "test@gmail.com"
$ws.Range("D4").Value = "Test123"
$ws.Range("E4").Value = "Male"
$ws.Range("F4").NumberFormat = "mm-dd-yy"
$ws.Range("F4").Value = [datetime]"1998-09-10"

# --- Row 5: unchanged (old "Test4") --------------------------------------
$ws.Range("A5").Value = "Test4"
$ws.Range("B5").Value = "j"
$ws.Range("C5").Value = "k"
$ws.Range("D5").Value = "l"

# --- Row 6: unchanged (old "Test5") --------------------------------------
$ws.Range("A6").Value = "Test5"
$ws.Range("B6").Value = "m"
$ws.Range("C6").Value = "n"
$ws.Range("D6").Value = "o"

# --- Row 7: unchanged (old "Test6") --------------------------------------
$ws.Range("A7").Value = "Test6"
$ws.Range("B7").Value = "p"
$ws.Range("C7").Value = "q"
$ws.Range("D7").Value = "r"

# Selection moves to F12 (below/right of the now-wider table)
[void]$ws.Range("F12").Select()
